# Update "想去人数" (F) and "最低票价" (G) figures across sheets to match
# the refreshed bilibili 漫展 listing data (gh-pages output regeneration).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 7955
$ws.Range("F5").Value = 2156
$ws.Range("F6").Value = 8686
$ws.Range("F8").Value = 98
$ws.Range("F9").Value = 5817
$ws.Range("F11").Value = 2841
$ws.Range("F12").Value = 1230
$ws.Range("F13").Value = 424
$ws.Range("F15").Value = 48
$ws.Range("F16").Value = 650
$ws.Range("G16").Value = 98
$ws.Range("F17").Value = 132
$ws.Range("F18").Value = 4075
$ws.Range("F20").Value = 74
$ws.Range("F21").Value = 77
$ws.Range("F24").Value = 38
$ws.Range("F25").Value = 5968
$ws.Range("F26").Value = 217
$ws.Range("F27").Value = 82
$ws.Range("F29").Value = 408
$ws.Range("F30").Value = 180
$ws.Range("F31").Value = 430
$ws.Range("F32").Value = 4302
$ws.Range("F33").Value = 1563
$ws.Range("F35").Value = 1729
$ws.Range("F36").Value = 5766
$ws.Range("F37").Value = 87
$ws.Range("F39").Value = 63
$ws.Range("F40").Value = 3774
$ws.Range("F41").Value = 39
$ws.Range("F42").Value = 57
$ws.Range("F43").Value = 11
$ws.Range("F44").Value = 2368
$ws.Range("F46").Value = 41
$ws.Range("F48").Value = 11
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 172
$ws.Range("G3").Value = 168
$ws.Range("F5").Value = 81
$ws.Range("F6").Value = 24
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1384
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1384
$ws.Range("F4").Value = 7955
$ws.Range("F5").Value = 2156
$ws.Range("F6").Value = 8686
$ws.Range("F7").Value = 98
$ws.Range("F8").Value = 5817
$ws.Range("F10").Value = 2841
$ws.Range("F11").Value = 1230
$ws.Range("F12").Value = 424
$ws.Range("F14").Value = 48
$ws.Range("F15").Value = 172
$ws.Range("G15").Value = 168
$ws.Range("F16").Value = 650
$ws.Range("G16").Value = 98
$ws.Range("F18").Value = 132
$ws.Range("F19").Value = 4075
$ws.Range("F21").Value = 74
$ws.Range("F22").Value = 77
$ws.Range("F25").Value = 38
$ws.Range("F26").Value = 5968
$ws.Range("F27").Value = 217
$ws.Range("F28").Value = 82
$ws.Range("F29").Value = 408
$ws.Range("F30").Value = 180
$ws.Range("F31").Value = 430
$ws.Range("F32").Value = 81
$ws.Range("F33").Value = 4302
$ws.Range("F34").Value = 1563
$ws.Range("F35").Value = 24
$ws.Range("F36").Value = 1729
$ws.Range("F38").Value = 5766
$ws.Range("F39").Value = 87
$ws.Range("F41").Value = 63
$ws.Range("F42").Value = 3774
$ws.Range("F43").Value = 57
$ws.Range("F44").Value = 2368
$ws.Range("F45").Value = 41
